$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a bare "NN%" string: Excel's smart-parsing would
# otherwise coerce them into a percentage NUMBER (e.g. 0.86) and rewrite the
# cell style. Force literal text by pre-formatting as Text, then restore the
# original (bordered, General-format) look so the style index is unchanged.
function Set-LiteralText($cell, $text) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
    $r.Borders.LineStyle = 1
}

$ws.Range("E2").Value = "2026-02-18 06:48:55"
$ws.Range("O2").Value = "-1.1 °C"
$ws.Range("E3").Value = "2026-02-18 06:48:58"
$ws.Range("O3").Value = "-3.6 °C"
$ws.Range("E4").Value = "2026-02-18 06:49:00"
$ws.Range("J4").Value = "1018.2 hPa"
$ws.Range("N4").Value = "4.2 °C 6:24 TU"
$ws.Range("O4").Value = "6.9 °C"
$ws.Range("E5").Value = "2026-02-18 06:49:03"
Set-LiteralText "H5" "86%"
$ws.Range("M5").Value = "0.3 °C 6:21 TU"
$ws.Range("O5").Value = "-1.8 °C"
$ws.Range("E6").Value = "2026-02-18 06:49:05"
$ws.Range("J6").Value = "1018.0 hPa"
$ws.Range("N6").Value = "6.0 °C 6:23 TU"
$ws.Range("O6").Value = "7.8 °C"
$ws.Range("E7").Value = "2026-02-18 06:49:08"
Set-LiteralText "H7" "83%"
$ws.Range("J7").Value = "1018.2 hPa"
$ws.Range("N7").Value = "11.0 °C 6:06 TU"
$ws.Range("E8").Value = "2026-02-18 06:49:10"
$ws.Range("J8").Value = "1018.3 hPa"
$ws.Range("O8").Value = "8.7 °C"
$ws.Range("E9").Value = "2026-02-18 06:49:13"
$ws.Range("N9").Value = "4.0 °C 6:00 TU"
$ws.Range("E10").Value = "2026-02-18 06:49:16"
Set-LiteralText "H10" "97%"
$ws.Range("L10").Value = "6.1 km/h - 76º 6:14 TU"
$ws.Range("N10").Value = "4.2 °C 6:02 TU"
$ws.Range("O10").Value = "6.6 °C"
$ws.Range("E11").Value = "2026-02-18 06:49:18"
$ws.Range("N11").Value = "-0.3 °C 6:26 TU"
$ws.Range("O11").Value = "1.8 °C"
$ws.Range("E12").Value = "2026-02-18 06:49:21"
$ws.Range("O12").Value = "5.8 °C"
$ws.Range("E13").Value = "2026-02-18 06:49:23"
$ws.Range("N13").Value = "-3.8 °C 6:00 TU"
$ws.Range("O13").Value = "-2.4 °C"
$ws.Range("E14").Value = "2026-02-18 06:49:26"
$ws.Range("E15").Value = "2026-02-18 06:49:28"
$ws.Range("E16").Value = "2026-02-18 06:49:31"
$ws.Range("E17").Value = "2026-02-18 06:49:33"
Set-LiteralText "H17" "86%"
$ws.Range("K17").Value = "-0.1 MJ/m2"
$ws.Range("O17").Value = "2.3 °C"
$ws.Range("E18").Value = "2026-02-18 06:49:36"
Set-LiteralText "H18" "96%"
$ws.Range("J18").Value = "1018.3 hPa"
$ws.Range("O18").Value = "7.2 °C"
$ws.Range("E19").Value = "2026-02-18 06:49:38"
$ws.Range("N19").Value = "5.2 °C 6:04 TU"
$ws.Range("E20").Value = "2026-02-18 06:49:41"
Set-LiteralText "H20" "77%"
$ws.Range("M20").Value = "-0.1 °C 6:10 TU"
$ws.Range("O20").Value = "-1.0 °C"
$ws.Range("E21").Value = "2026-02-18 06:49:43"
Set-LiteralText "H21" "85%"
$ws.Range("N21").Value = "0.2 °C 6:03 TU"
$ws.Range("O21").Value = "1.7 °C"
$ws.Range("E22").Value = "2026-02-18 06:49:46"
$ws.Range("E23").Value = "2026-02-18 06:49:48"
$ws.Range("O23").Value = "0.6 °C"
$ws.Range("E24").Value = "2026-02-18 06:49:51"
$ws.Range("J24").Value = "1018.8 hPa"
$ws.Range("N24").Value = "2.4 °C 6:03 TU"
$ws.Range("O24").Value = "4.9 °C"
$ws.Range("E25").Value = "2026-02-18 06:49:53"
$ws.Range("M25").Value = "1.5 °C 6:26 TU"
$ws.Range("O25").Value = "-0.2 °C"
$ws.Range("E26").Value = "2026-02-18 06:49:56"
$ws.Range("E27").Value = "2026-02-18 06:49:58"
Set-LiteralText "H27" "50%"
$ws.Range("K27").Value = "-0.1 MJ/m2"
$ws.Range("O27").Value = "0.9 °C"
$ws.Range("E28").Value = "2026-02-18 06:50:01"
$ws.Range("J28").Value = "1018.7 hPa"
$ws.Range("O28").Value = "4.6 °C"
$ws.Range("E29").Value = "2026-02-18 06:50:03"
$ws.Range("E30").Value = "2026-02-18 06:50:06"
$ws.Range("J30").Value = "1018.3 hPa"
$ws.Range("E31").Value = "2026-02-18 06:50:08"
$ws.Range("J31").Value = "1016.9 hPa"
$ws.Range("E32").Value = "2026-02-18 06:50:11"
$ws.Range("E33").Value = "2026-02-18 06:50:14"
$ws.Range("J33").Value = "1021.2 hPa"
$ws.Range("O33").Value = "-0.2 °C"
$ws.Range("E34").Value = "2026-02-18 06:50:16"
Set-LiteralText "H34" "52%"
$ws.Range("M34").Value = "3.7 °C 6:23 TU"
$ws.Range("O34").Value = "0.4 °C"
$ws.Range("E35").Value = "2026-02-18 06:50:19"
$ws.Range("N35").Value = "4.2 °C 6:22 TU"
$ws.Range("O35").Value = "7.0 °C"
$ws.Range("E36").Value = "2026-02-18 06:50:21"
$ws.Range("E37").Value = "2026-02-18 06:50:24"
$ws.Range("J37").Value = "1021.1 hPa"
$ws.Range("N37").Value = "-0.2 °C 6:08 TU"
$ws.Range("O37").Value = "1.2 °C"
$ws.Range("E38").Value = "2026-02-18 06:50:26"
Set-LiteralText "H38" "94%"
$ws.Range("N38").Value = "5.8 °C 6:27 TU"
$ws.Range("O38").Value = "8.5 °C"
$ws.Range("E39").Value = "2026-02-18 06:50:29"
$ws.Range("O39").Value = "0.4 °C"
$ws.Range("E40").Value = "2026-02-18 06:50:31"
$ws.Range("N40").Value = "-0.8 °C 6:29 TU"
$ws.Range("O40").Value = "0.6 °C"
$ws.Range("E41").Value = "2026-02-18 06:50:34"
Set-LiteralText "H41" "94%"
$ws.Range("J41").Value = "1017.9 hPa"
$ws.Range("O41").Value = "8.2 °C"
$ws.Range("E42").Value = "2026-02-18 06:50:36"
$ws.Range("O42").Value = "7.9 °C"
$ws.Range("E43").Value = "2026-02-18 06:50:39"
$ws.Range("N43").Value = "5.1 °C 6:29 TU"
$ws.Range("O43").Value = "6.9 °C"
$ws.Range("E44").Value = "2026-02-18 06:50:41"
Set-LiteralText "H44" "64%"
$ws.Range("E45").Value = "2026-02-18 06:50:43"
Set-LiteralText "H45" "95%"
$ws.Range("J45").Value = "1020.3 hPa"
$ws.Range("E46").Value = "2026-02-18 06:50:46"
$ws.Range("N46").Value = "4.2 °C 6:03 TU"
$ws.Range("O46").Value = "6.0 °C"
